$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated PCA/LDA transformed coordinates for columns A (PC1) and B (PC2).
# Column C (Labels) is unchanged.
$ws.Cells.Item(2, 1).Value = 22.58364190325597
$ws.Cells.Item(2, 2).Value = -56.13496713645615
$ws.Cells.Item(3, 1).Value = 43.60533359809029
$ws.Cells.Item(3, 2).Value = -24.37124529499845
$ws.Cells.Item(4, 1).Value = 31.77247806104575
$ws.Cells.Item(4, 2).Value = 43.89449190627568
$ws.Cells.Item(5, 1).Value = 36.22098572794811
$ws.Cells.Item(5, 2).Value = 61.79061866308648
$ws.Cells.Item(6, 1).Value = 38.43087430308902
$ws.Cells.Item(6, 2).Value = 65.39004204615279
$ws.Cells.Item(7, 1).Value = 40.38613145275494
$ws.Cells.Item(7, 2).Value = 69.56176818315582
$ws.Cells.Item(8, 1).Value = 42.00579826283618
$ws.Cells.Item(8, 2).Value = 76.68794164859358
$ws.Cells.Item(9, 1).Value = 52.19136477635872
$ws.Cells.Item(9, 2).Value = 99.25097515595091
$ws.Cells.Item(10, 1).Value = 68.88199345259925
$ws.Cells.Item(10, 2).Value = 118.948943988071
$ws.Cells.Item(11, 1).Value = 96.86753007657389
$ws.Cells.Item(11, 2).Value = 135.4551514497722
$ws.Cells.Item(12, 1).Value = 34.12364113983868
$ws.Cells.Item(12, 2).Value = 30.745165234074
$ws.Cells.Item(13, 1).Value = 44.04136021147245
$ws.Cells.Item(13, 2).Value = -24.48508649461825
$ws.Cells.Item(14, 1).Value = 22.56625043833582
$ws.Cells.Item(14, 2).Value = -55.61746398950658
$ws.Cells.Item(15, 1).Value = 33.87625747777388
$ws.Cells.Item(15, 2).Value = 30.70664538771633
$ws.Cells.Item(16, 1).Value = 43.18716545869167
$ws.Cells.Item(16, 2).Value = -24.15976399227024
$ws.Cells.Item(17, 1).Value = 22.74602341084481
$ws.Cells.Item(17, 2).Value = -57.30453828288847
$ws.Cells.Item(18, 1).Value = 34.30377328284699
$ws.Cells.Item(18, 2).Value = 30.77321341807631
$ws.Cells.Item(19, 1).Value = 43.44151134737073
$ws.Cells.Item(19, 2).Value = -23.78409461006696
$ws.Cells.Item(20, 1).Value = 22.24107886174357
$ws.Cells.Item(20, 2).Value = -55.74834558300439
$ws.Cells.Item(21, 1).Value = 33.97510633590857
$ws.Cells.Item(21, 2).Value = 30.722037038063
$ws.Cells.Item(22, 1).Value = 43.90576614420564
$ws.Cells.Item(22, 2).Value = -23.114593173397
$ws.Cells.Item(23, 1).Value = 22.6401763833145
$ws.Cells.Item(23, 2).Value = -56.00895921534407
$ws.Cells.Item(24, 1).Value = 33.72236509655497
$ws.Cells.Item(24, 2).Value = 30.68268296905581
$ws.Cells.Item(25, 1).Value = 42.88575725458053
$ws.Cells.Item(25, 2).Value = -24.50814269726197
$ws.Cells.Item(26, 1).Value = 22.83990860370458
$ws.Cells.Item(26, 2).Value = -56.48934482421513
$ws.Cells.Item(27, 1).Value = 46.64464959041608
$ws.Cells.Item(27, 2).Value = -77.27232414422173
$ws.Cells.Item(28, 1).Value = 204.6111235130482
$ws.Cells.Item(28, 2).Value = -34.34135917735852
$ws.Cells.Item(29, 1).Value = 242.791864001831
$ws.Cells.Item(29, 2).Value = -5.399583334040813
$ws.Cells.Item(30, 1).Value = 207.5071627221974
$ws.Cells.Item(30, 2).Value = 53.79373917523524
$ws.Cells.Item(31, 1).Value = 191.0452012094882
$ws.Cells.Item(31, 2).Value = 67.2631125309501
$ws.Cells.Item(32, 1).Value = 181.2236385963581
$ws.Cells.Item(32, 2).Value = 69.91719393286874
$ws.Cells.Item(33, 1).Value = 184.73479978106
$ws.Cells.Item(33, 2).Value = 72.57512592501678
$ws.Cells.Item(34, 1).Value = 94.76692141648324
$ws.Cells.Item(34, 2).Value = 59.94085410569714
$ws.Cells.Item(35, 1).Value = 80.76858059347963
$ws.Cells.Item(35, 2).Value = 45.48417748526482
$ws.Cells.Item(36, 1).Value = 78.02263138591995
$ws.Cells.Item(36, 2).Value = 66.98925742493343
$ws.Cells.Item(37, 1).Value = 241.9944098317442
$ws.Cells.Item(37, 2).Value = -5.710334351748137
$ws.Cells.Item(38, 1).Value = 205.5465410488141
$ws.Cells.Item(38, 2).Value = -34.40770169725737
$ws.Cells.Item(39, 1).Value = 45.37154394486111
$ws.Cells.Item(39, 2).Value = -77.71824954124436
$ws.Cells.Item(40, 1).Value = 244.600656698791
$ws.Cells.Item(40, 2).Value = -5.343773879081229
$ws.Cells.Item(41, 1).Value = 202.4549850797038
$ws.Cells.Item(41, 2).Value = -34.52179637043669
$ws.Cells.Item(42, 1).Value = 46.24457562299619
$ws.Cells.Item(42, 2).Value = -77.73429930462642
$ws.Cells.Item(43, 1).Value = 243.2544212442867
$ws.Cells.Item(43, 2).Value = -5.038837443312665
$ws.Cells.Item(44, 1).Value = 202.8919953077074
$ws.Cells.Item(44, 2).Value = -34.43360737771537
$ws.Cells.Item(45, 1).Value = 48.25576137824801
$ws.Cells.Item(45, 2).Value = -77.02144300597661
$ws.Cells.Item(46, 1).Value = 243.6678471207606
$ws.Cells.Item(46, 2).Value = -5.497095720359951
$ws.Cells.Item(47, 1).Value = 203.7444567600446
$ws.Cells.Item(47, 2).Value = -34.47718257981345
$ws.Cells.Item(48, 1).Value = 47.68925121723031
$ws.Cells.Item(48, 2).Value = -77.02796429514441
$ws.Cells.Item(49, 1).Value = 244.753444092229
$ws.Cells.Item(49, 2).Value = -4.8186073134851
$ws.Cells.Item(50, 1).Value = 201.7702559551056
$ws.Cells.Item(50, 2).Value = -34.99790505699498
$ws.Cells.Item(51, 1).Value = 47.56465293619551
$ws.Cells.Item(51, 2).Value = -76.96093547347535
$ws.Cells.Item(52, 1).Value = -151.0940430931603
$ws.Cells.Item(52, 2).Value = -59.95028446452336
$ws.Cells.Item(53, 1).Value = -183.4403339199195
$ws.Cells.Item(53, 2).Value = -31.6416146917306
$ws.Cells.Item(54, 1).Value = -212.7012446593911
$ws.Cells.Item(54, 2).Value = 5.135791612638164
$ws.Cells.Item(55, 1).Value = -234.8783794805104
$ws.Cells.Item(55, 2).Value = 5.533155899488222
$ws.Cells.Item(56, 1).Value = -237.2371440164432
$ws.Cells.Item(56, 2).Value = 56.67243603083977
$ws.Cells.Item(57, 1).Value = -237.5868694833788
$ws.Cells.Item(57, 2).Value = 74.56048887874405
$ws.Cells.Item(58, 1).Value = -238.211076222787
$ws.Cells.Item(58, 2).Value = 73.52693939817121
$ws.Cells.Item(59, 1).Value = -238.5072291916436
$ws.Cells.Item(59, 2).Value = 73.65970298305446
$ws.Cells.Item(60, 1).Value = -240.208019071633
$ws.Cells.Item(60, 2).Value = 75.14042885451154
$ws.Cells.Item(61, 1).Value = -240.5605929758892
$ws.Cells.Item(61, 2).Value = 84.89379608371982
$ws.Cells.Item(62, 1).Value = -213.0764002147652
$ws.Cells.Item(62, 2).Value = 5.147134617352033
$ws.Cells.Item(63, 1).Value = -183.6352825227974
$ws.Cells.Item(63, 2).Value = -31.6686150162398
$ws.Cells.Item(64, 1).Value = -151.3514581631625
$ws.Cells.Item(64, 2).Value = -60.68661761880731
$ws.Cells.Item(65, 1).Value = -212.5119001053074
$ws.Cells.Item(65, 2).Value = 4.94567733319325
$ws.Cells.Item(66, 1).Value = -184.2558729954786
$ws.Cells.Item(66, 2).Value = -33.00452757324563
$ws.Cells.Item(67, 1).Value = -152.0504885054203
$ws.Cells.Item(67, 2).Value = -59.21444910971279
$ws.Cells.Item(68, 1).Value = -213.27543104319
$ws.Cells.Item(68, 2).Value = 4.800725859616044
$ws.Cells.Item(69, 1).Value = -184.1648396893998
$ws.Cells.Item(69, 2).Value = -32.1019478737936
$ws.Cells.Item(70, 1).Value = -150.6914893976012
$ws.Cells.Item(70, 2).Value = -59.06342363955449
$ws.Cells.Item(71, 1).Value = -213.0921113672618
$ws.Cells.Item(71, 2).Value = 5.310156941893615
$ws.Cells.Item(72, 1).Value = -183.8322746961003
$ws.Cells.Item(72, 2).Value = -31.29371860763054
$ws.Cells.Item(73, 1).Value = -150.9773576786609
$ws.Cells.Item(73, 2).Value = -59.69075689436248
$ws.Cells.Item(74, 1).Value = -212.8367125781967
$ws.Cells.Item(74, 2).Value = 5.167266813379273
$ws.Cells.Item(75, 1).Value = -184.0348571069301
$ws.Cells.Item(75, 2).Value = -32.10247155159171
$ws.Cells.Item(76, 1).Value = -151.1522359317104
$ws.Cells.Item(76, 2).Value = -60.1988665730984
